# Update assessment results: renumber userIds, add new students, clear
# scores to "NA" for everyone except the last row (Ishdeep Singh), who
# now carries numeric scores and moves to the bottom of the new list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# userId / firstName / lastName / D / E / F for rows 2..10
$data = @(
    @(10, "Yvan",     "Labiche",    "NA", "NA", "NA"),
    @(11, "Jenna",    "McConnell",  "NA", "NA", "NA"),
    @(12, "Jennifer", "Poll",       "NA", "NA", "NA"),
    @(13, "Saundra",  "Warmington", "NA", "NA", "NA"),
    @(14, "Erica",    "East",       "NA", "NA", "NA"),
    @(15, "Jerry ",   "Buburuz",    "NA", "NA", "NA"),
    @(16, "Daren ",   "Russ",       "NA", "NA", "NA"),
    @(17, "kong ",    "Chiv",       "NA", "NA", "NA"),
    @(18, "Ishdeep",  "Singh",      4,    4,    4)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
